# The grade threshold used when normalizing scores (stored in $C$29, the
# class average) is raised from 76 to 77. This formula lives in D2 (its own
# formula) and is the anchor of the shared formula spanning D3:D28, so both
# need to be updated; updating D3:D28 together keeps them sharing one formula
# group, matching the workbook's original layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=C2 + (IF(`$C`$29 < 77, 77, `$C`$29) - `$C`$29)"
$ws.Range("D3:D28").Formula = "=C3 + (IF(`$C`$29 < 77, 77, `$C`$29) - `$C`$29)"
